$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Metadata" sheet: update URL, Version, Date, Publisher
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-response-item-status"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------------------
# "Elements" sheet: clear top-level Extension row's Constraint(s) cell,
# widen column Y, and update the (now renamed) payer-claim-status ValueSet url
# ---------------------------------------------------------------------------
$el = $wb.Worksheets.Item("Elements")
$el.Range("AI2").Value = ""
$el.Columns.Item(25).ColumnWidth = 57.8671875
$el.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/payer-claim-status"
